$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for the latest weekly record, pushing the existing
# history down by one row (row 2 — the most-recent week already on file —
# stays put; everything from the old row 3 onward shifts to row+1).
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with this week's record.
$ws.Cells.Item(3, 1).Value = 11
$ws.Cells.Item(3, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(3, 3).Value = "Bíobío"
$ws.Cells.Item(3, 4).Value = 44462
$ws.Cells.Item(3, 5).Value = 8
$ws.Cells.Item(3, 6).Value = 100114007
$ws.Cells.Item(3, 7).Value = "Jengibre"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 60
$ws.Cells.Item(3, 11).Value = 14000
$ws.Cells.Item(3, 12).Value = 15000
$ws.Cells.Item(3, 13).Value = 14500
$ws.Cells.Item(3, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(3, 15).Value = "Perú"
$ws.Cells.Item(3, 16).Value = 1115
$ws.Cells.Item(3, 17).Value = 13
$ws.Cells.Item(3, 18).Value = "Hortaliza"
